$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("getpages_with_index")

# Step 1: insert a new row at position 3 (shifts old rows 3,4 down to 4,5)
$ws1.Rows.Item(3).Insert() | Out-Null
$ws1.Cells.Item(3, 1).Value = "/espanol/efectos-secundarios"
$ws1.Cells.Item(3, 2).Value = "Article"
$ws1.Cells.Item(3, 3).Value = "Spanish"

# Step 2: append new rows 6 through 24
$ws1.Cells.Item(6, 1).Value = "/espanol/noticias/comunicados-de-prensa/2018/leucemia-llc-ibrutinib-estudio"
$ws1.Cells.Item(6, 2).Value = "Press Release"
$ws1.Cells.Item(6, 3).Value = "Spanish"
$ws1.Cells.Item(7, 1).Value = "/news-events/press-releases/2018/leukemia-cll-ibrutinib-trial"
$ws1.Cells.Item(7, 2).Value = "Press Release"
$ws1.Cells.Item(7, 3).Value = "English"
$ws1.Cells.Item(8, 1).Value = "/news-events/cancer-currents-blog"
$ws1.Cells.Item(8, 2).Value = "Blog Series"
$ws1.Cells.Item(8, 3).Value = "English"
$ws1.Cells.Item(9, 1).Value = "/espanol/noticias/temas-y-relatos-blog"
$ws1.Cells.Item(9, 2).Value = "Blog Series"
$ws1.Cells.Item(9, 3).Value = "Spanish"
$ws1.Cells.Item(10, 1).Value = "/about-cancer/coping/feelings/relaxation/chanock-stephen"
$ws1.Cells.Item(10, 2).Value = "Biography"
$ws1.Cells.Item(10, 3).Value = "English"
$ws1.Cells.Item(11, 1).Value = "/about-cancer/coping/feelings/relaxation/dfharvard"
$ws1.Cells.Item(11, 2).Value = "Cancer Center"
$ws1.Cells.Item(11, 3).Value = "English"
$ws1.Cells.Item(12, 1).Value = "/types/breast/research"
$ws1.Cells.Item(12, 2).Value = "Cancer Research List Page"
$ws1.Cells.Item(12, 3).Value = "English"
$ws1.Cells.Item(13, 1).Value = "/espanol/tipos/seno/investigacion"
$ws1.Cells.Item(13, 2).Value = "Cancer Research List Page"
$ws1.Cells.Item(13, 3).Value = "Spanish"
$ws1.Cells.Item(14, 1).Value = "/types/breast"
$ws1.Cells.Item(14, 2).Value = "Cancer Type Home Page"
$ws1.Cells.Item(14, 3).Value = "English"
$ws1.Cells.Item(15, 1).Value = "/espanol/tipos/seno"
$ws1.Cells.Item(15, 2).Value = "Cancer Type Home Page"
$ws1.Cells.Item(15, 3).Value = "Spanish"
$ws1.Cells.Item(16, 1).Value = "/types/breast/hp"
$ws1.Cells.Item(16, 2).Value = "Cancer Type Home Page"
$ws1.Cells.Item(16, 3).Value = "English"
$ws1.Cells.Item(17, 1).Value = "/about-cancer/coping/feelings/relaxation/2019-investigators-site"
$ws1.Cells.Item(17, 2).Value = "Event"
$ws1.Cells.Item(17, 3).Value = "English"
$ws1.Cells.Item(18, 1).Value = "/about-cancer"
$ws1.Cells.Item(18, 2).Value = "Home & Landing"
$ws1.Cells.Item(18, 3).Value = "English"
$ws1.Cells.Item(19, 1).Value = "/espanol/cancer"
$ws1.Cells.Item(19, 2).Value = "Home & Landing"
$ws1.Cells.Item(19, 3).Value = "Spanish"
$ws1.Cells.Item(20, 1).Value = "/news-events/press-releases/2018"
$ws1.Cells.Item(20, 2).Value = "Mini Landing Page"
$ws1.Cells.Item(20, 3).Value = "English"
$ws1.Cells.Item(21, 1).Value = "/espanol/noticias/comunicados-de-prensa/2018"
$ws1.Cells.Item(21, 2).Value = "Mini Landing Page"
$ws1.Cells.Item(21, 3).Value = "Spanish"
$ws1.Cells.Item(22, 1).Value = "/about-nci/organization/screen-to-save-infographic"
$ws1.Cells.Item(22, 2).Value = "Infographic"
$ws1.Cells.Item(22, 3).Value = "English"
$ws1.Cells.Item(23, 1).Value = "/espanol/infografia-nci"
$ws1.Cells.Item(23, 2).Value = "Infographic"
$ws1.Cells.Item(23, 3).Value = "Spanish"
$ws1.Cells.Item(24, 1).Value = "/research/progress/discovery/gutcheck-intro-video"
$ws1.Cells.Item(24, 2).Value = "Video"
$ws1.Cells.Item(24, 3).Value = "English"

# Step 3: update existing cell B2 ("article" -> "Article") last
$ws1.Cells.Item(2, 2).Value = "Article"

# Activate sheet1 and select I21 to match target sheetView/selection
$ws1.Activate() | Out-Null
$ws1.Range("I21").Select() | Out-Null
